$d = $word.ActiveDocument

# 1. Title: merge "Product " + "Backlog :" runs into "Product Backlog :"
$d.Content.Find.Execute("Product Backlog :", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Product Backlog :", 2) | Out-Null

# 2. Title: merge "Whatsapp" + " " into "Whatsapp "
$d.Content.Find.Execute("Whatsapp ", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Whatsapp ", 2) | Out-Null

# 3. Merge "See the time and date of a text in a group when it was " + "send"
$d.Content.Find.Execute("See the time and date of a text in a group when it was send",
                         $false, $false, $false, $false, $false,
                         $true, 1, $false, "See the time and date of a text in a group when it was send", 2) | Out-Null

# 4. Status changes: first two "In Progress" -> "Done", then two "To be started" -> "Done"
$d.Content.Find.Execute("In Progress", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Done", 2) | Out-Null
$d.Content.Find.Execute("In Progress", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Done", 2) | Out-Null
$d.Content.Find.Execute("To be started", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Done", 2) | Out-Null
$d.Content.Find.Execute("To be started", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Done", 2) | Out-Null
